# Updated cryptos list with GitHub Actions (matches the scraper's latest run).
# Price (D) / Volume(1h) (E) refreshed for each coin row; MXToken and
# RenderToken swapped rank positions (rows 40/41).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A handful of Price cells carry significant trailing zeros (e.g. "1.0000")
# that Excel would normally collapse if the string is auto-recognised as a
# number. Prefixing with a literal leading apostrophe ('' inside a
# single-quoted PowerShell string escapes to a literal ' character) forces
# Excel to keep the value as text, exactly like typing '1.0000 in the UI.

$ws.Range('D2').Value = '26.600.93'
$ws.Range('E2').Value = '  +4.14%  '
$ws.Range('D3').Value = '1.743.04'
$ws.Range('E3').Value = '  +4.41%  '
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').Value = '246.52'
$ws.Range('E5').Value = '  +4.26%  '
$ws.Range('D6').Value = '''1.0000'
$ws.Range('E6').Value = '  -0.03%  '
$ws.Range('D7').Value = '0.4823'
$ws.Range('E7').Value = '  +2.02%  '
$ws.Range('D8').Value = '0.2693'
$ws.Range('E8').Value = '  +3.82%  '
$ws.Range('D9').Value = '0.06267'
$ws.Range('E9').Value = '  +1.74%  '
$ws.Range('D10').Value = '1.744.61'
$ws.Range('E10').Value = '  +4.60%  '
$ws.Range('D11').Value = '0.07128'
$ws.Range('E11').Value = '  +1.87%  '
$ws.Range('D12').Value = '''15.90'
$ws.Range('E12').Value = '  +7.71%  '
$ws.Range('E13').Value = '  +6.93%  '
$ws.Range('D14').Value = '4.517'
$ws.Range('E14').Value = '  +3.58%  '
$ws.Range('D15').Value = '77.44'
$ws.Range('E15').Value = '  +2.71%  '
$ws.Range('D17').Value = '26.603.79'
$ws.Range('E17').Value = '  +4.16%  '
$ws.Range('D18').Value = '0.9999'
$ws.Range('E18').Value = '  -0.08%  '
$ws.Range('D19').Value = '11.82'
$ws.Range('E19').Value = '  +3.68%  '
$ws.Range('E20').Value = '  +2.74%  '
$ws.Range('D21').Value = '1.969.48'
$ws.Range('E21').Value = '  +4.62%  '
$ws.Range('D22').Value = '4.618'
$ws.Range('E22').Value = '  +4.14%  '
$ws.Range('D23').Value = '8.882'
$ws.Range('E23').Value = '  +1.36%  '
$ws.Range('D24').Value = '5.365'
$ws.Range('E24').Value = '  +2.62%  '
$ws.Range('D25').Value = '''136.50'
$ws.Range('E25').Value = '  -0.38%  '
$ws.Range('D26').Value = '15.37'
$ws.Range('E26').Value = '  +2.72%  '
$ws.Range('E27').Value = '  +5.87%  '
$ws.Range('E28').Value = '  +3.09%  '
$ws.Range('D29').Value = '106.78'
$ws.Range('E29').Value = '  +2.26%  '
$ws.Range('D30').Value = '''4.020'
$ws.Range('E30').Value = '  +0.54%  '
$ws.Range('D31').Value = '3.741'
$ws.Range('E31').Value = '  +3.18%  '
$ws.Range('D32').Value = '0.07891'
$ws.Range('E32').Value = '  +0.90%  '
$ws.Range('D33').Value = '''0.04590'
$ws.Range('E33').Value = '  +6.94%  '
$ws.Range('D34').Value = '0.9996'
$ws.Range('E34').Value = '  +0.03%  '
$ws.Range('D36').Value = '0.6409'
$ws.Range('E36').Value = '  +5.93%  '
$ws.Range('D37').Value = '1.001'
$ws.Range('E37').Value = '  +5.01%  '
$ws.Range('D38').Value = '0.9317'
$ws.Range('E38').Value = '  -0.06%  '
$ws.Range('D39').Value = '113.86'
$ws.Range('E39').Value = '  +14.12%  '
$ws.Range('B40').Value = 'RenderToken'
$ws.Range('C40').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D40').Value = '1.996'
$ws.Range('E40').Value = '  +8.18%  '
$ws.Range('B41').Value = 'MXToken'
$ws.Range('C41').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D41').Value = '2.444'
$ws.Range('E41').Value = '  -3.13%  '
$ws.Range('D42').Value = '5.805'
$ws.Range('E42').Value = '  +18.67%  '
$ws.Range('E43').Value = '  +0.01%  '
$ws.Range('D44').Value = '0.01511'
$ws.Range('E44').Value = '  +2.19%  '
$ws.Range('D45').Value = '''0.3920'
$ws.Range('E45').Value = '  +4.78%  '
$ws.Range('D46').Value = '0.1217'
$ws.Range('E46').Value = '  +9.27%  '
$ws.Range('D47').Value = '''6.750'
$ws.Range('E47').Value = '  +8.99%  '
$ws.Range('E48').Value = '  +1.41%  '
$ws.Range('D49').Value = '7.976'
$ws.Range('E49').Value = '  +7.29%  '
$ws.Range('D50').Value = '30.79'
$ws.Range('E50').Value = '  +3.34%  '
$ws.Range('D51').Value = '1.264'
$ws.Range('E51').Value = '  +5.15%  '
